$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: Status column (B) and Latest Handoff Date column (C) both
# previously held "Ready for handoff" text for both file rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# New "Latest Target File" (F) / "Latest Handback File" (G) hyperlinks, mirroring
# the existing "Latest Handoff File" (A) / target xlf (D) hyperlinks for each row.
$bdffAddr = $wsZh.Hyperlinks.Item(1).Address
$bdffXlfAddr = $wsZh.Hyperlinks.Item(2).Address
$fdAddr = $wsZh.Hyperlinks.Item(3).Address
$fdXlfAddr = $wsZh.Hyperlinks.Item(4).Address

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $bdffAddr, "", "", "bdff6a9b-133e-4bf0-8112-ba6ecb2eaadb.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $bdffXlfAddr, "", "", "bdff6a9b-133e-4bf0-8112-ba6ecb2eaadb.841fb03b39a3fd9740e3398c2b9e17768485749a.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $fdAddr, "", "", "fd0d4ac0-7a18-42c1-b385-ec68358957ff.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $fdXlfAddr, "", "", "fd0d4ac0-7a18-42c1-b385-ec68358957ff.2b3c79a37eb42242be479354355d52e5a71d2274.zh-cn.xlf")

$wsZh.Range("F2").Style = "HyperLink"
$wsZh.Range("G2").Style = "HyperLink"
$wsZh.Range("F3").Style = "HyperLink"
$wsZh.Range("G3").Style = "HyperLink"

# Latest Handback DateTime column (H) now populated with the real handback time
$wsZh.Range("H2").Value = "2016-03-21 18:14:49"
$wsZh.Range("H3").Value = "2016-03-21 18:14:49"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$bdffAddrDe = $wsDe.Hyperlinks.Item(1).Address
$bdffXlfAddrDe = $wsDe.Hyperlinks.Item(2).Address
$fdAddrDe = $wsDe.Hyperlinks.Item(3).Address
$fdXlfAddrDe = $wsDe.Hyperlinks.Item(4).Address

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $bdffAddrDe, "", "", "bdff6a9b-133e-4bf0-8112-ba6ecb2eaadb.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $bdffXlfAddrDe, "", "", "bdff6a9b-133e-4bf0-8112-ba6ecb2eaadb.841fb03b39a3fd9740e3398c2b9e17768485749a.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $fdAddrDe, "", "", "fd0d4ac0-7a18-42c1-b385-ec68358957ff.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $fdXlfAddrDe, "", "", "fd0d4ac0-7a18-42c1-b385-ec68358957ff.2b3c79a37eb42242be479354355d52e5a71d2274.de-de.xlf")

$wsDe.Range("F2").Style = "HyperLink"
$wsDe.Range("G2").Style = "HyperLink"
$wsDe.Range("F3").Style = "HyperLink"
$wsDe.Range("G3").Style = "HyperLink"

# Latest Handback DateTime column (H) now populated with the real handback time
$wsDe.Range("H2").Value = "2016-03-21 18:14:55"
$wsDe.Range("H3").Value = "2016-03-21 18:14:55"
